$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Additional Sheet")

# Shift the existing content down and insert a new header row
$ws.Range("A1").Value = "This is the second sheet"
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Text"

$ws.Range("A2").Select()
